$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (rich shared strings) - Volume/Number + report week
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/17/2023  Through  7/23/2023"

# ---------------------------------------------------------------------------
# Helper number formats matching the sheet's existing styles
#   style 15 -> "#,##0"                      (plain integer counts)
#   style 16 -> '#,##0.0;"-"#,##0.0'          (percent-change figures)
# ---------------------------------------------------------------------------
$fmtInt = "#,##0"
$fmtPct = '#,##0.0;"-"#,##0.0'

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
$ws.Range("D14").Copy($ws.Range("C15"))   # numeric -> text "0" (reuse existing text cell's type/style)
$ws.Range("D15").NumberFormat = $fmtInt
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = $fmtPct
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = -16.666666666666
$ws.Range("L15").Value = -9.090909090909

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 14.285714285714
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = 22.222222222222
$ws.Range("I16").Value = 94
$ws.Range("J16").Value = 115
$ws.Range("K16").Value = -18.260869565217
$ws.Range("L16").Value = 54.098360655737
$ws.Range("M16").Value = -37.333333333333
$ws.Range("N16").Value = -84.887459807074

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 35
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 207
$ws.Range("J17").Value = 224
$ws.Range("K17").Value = -7.589285714285
$ws.Range("L17").Value = 40.816326530612
$ws.Range("M17").Value = 68.292682926829
$ws.Range("N17").Value = 5.612244897959

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -16.666666666666
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -7.142857142857
$ws.Range("I18").Value = 87
$ws.Range("J18").Value = 78
$ws.Range("K18").Value = 11.538461538461
$ws.Range("L18").Value = 33.846153846153
$ws.Range("M18").Value = -48.823529411764
$ws.Range("N18").Value = -88.671875

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 28.571428571428
$ws.Range("F19").Value = 27
$ws.Range("H19").Value = 8
$ws.Range("I19").Value = 176
$ws.Range("J19").Value = 197
$ws.Range("K19").Value = -10.659898477157
$ws.Range("L19").Value = 45.454545454545
$ws.Range("M19").Value = -5.376344086021
$ws.Range("N19").Value = -41.333333333333

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 11
$ws.Range("E20").Value = 83.333333333333
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = 13.636363636363
$ws.Range("I20").Value = 153
$ws.Range("J20").Value = 139
$ws.Range("K20").Value = 10.071942446043
$ws.Range("L20").Value = 25.409836065573
$ws.Range("M20").Value = 6.993006993006
$ws.Range("N20").Value = -92.080745341614

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = -5.128205128205
$ws.Range("F21").Value = 125
$ws.Range("G21").Value = 127
$ws.Range("H21").Value = -1.574803149606
$ws.Range("I21").Value = 731
$ws.Range("J21").Value = 768
$ws.Range("K21").Value = -4.817708333333
$ws.Range("L21").Value = 37.924528301886
$ws.Range("M21").Value = -6.402048655569
$ws.Range("N21").Value = -80.968497787034

# ---------------------------------------------------------------------------
# Row 22 - Transit (was all dashes/text, now has real figures)
# ---------------------------------------------------------------------------
$ws.Range("C22").NumberFormat = $fmtInt
$ws.Range("C22").Value = 1
$ws.Range("D22").NumberFormat = $fmtInt
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = $fmtPct
$ws.Range("E22").Value = 0
$ws.Range("F22").NumberFormat = $fmtInt
$ws.Range("F22").Value = 1
$ws.Range("G22").NumberFormat = $fmtInt
$ws.Range("G22").Value = 1
$ws.Range("H22").NumberFormat = $fmtPct
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 6
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = -40
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = -33.333333333333

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -10.714285714285
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 127
$ws.Range("H24").Value = -24.409448818897
$ws.Range("I24").Value = 837
$ws.Range("J24").Value = 851
$ws.Range("K24").Value = -1.645123384253
$ws.Range("L24").Value = 53.016453382084
$ws.Range("M24").Value = 88.089887640449

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 13
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 54
$ws.Range("H25").Value = -38.888888888888
$ws.Range("I25").Value = 302
$ws.Range("J25").Value = 305
$ws.Range("K25").Value = -0.983606557377
$ws.Range("L25").Value = 21.774193548387
$ws.Range("M25").Value = -14.204545454545

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------------
$ws.Range("D14").Copy($ws.Range("C26"))   # numeric -> text "0"
$ws.Range("D26").NumberFormat = $fmtInt
$ws.Range("D26").Value = 1
$ws.Range("E26").NumberFormat = $fmtPct
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = -66.666666666666
$ws.Range("J26").Value = 23
$ws.Range("K26").Value = -13.043478260869
$ws.Range("L26").Value = 5.263157894736

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("C27").Value = 1
$ws.Range("D27").NumberFormat = $fmtInt
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = $fmtPct
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = 150
$ws.Range("I27").Value = 38
$ws.Range("J27").Value = 26
$ws.Range("K27").Value = 46.153846153846
$ws.Range("L27").Value = -2.564102564102
